# Applies the "automatic update of files" edit:
#  1. Column C ("Förändrad") for rows 2..97 changes from 45184 to 45186.
#  2. HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2..12 gain a
#     second argument (the friendly name), equal to the value in column A
#     of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 97
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    # Update the "Förändrad" date in column C.
    $ws.Range("C$r").Value = 45186

    # Append the friendly-name argument to any HYPERLINK formulas on this row.
    $name = $ws.Range("A$r").Value()

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula()
        if ($f -and $f.Length -gt 0 -and $f.ToUpper().Contains("HYPERLINK(") -and -not $f.Contains(",")) {
            $newF = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
            $cell.Formula = $newF
        }
    }
}
